$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = "system, backup@backdoor.com, System"
    3 = "System, dnasr281@gmail.com"
    4 = "backup@backdoor.com, System"
    5 = "backup@backdoor.com, System"
    6 = "System, dnasr281@gmail.com"
    8 = "backup@backdoor.com, System"
    10 = "System, dnasr281@gmail.com"
    11 = "System, dnasr281@gmail.com"
    12 = "System, dnasr281@gmail.com"
    13 = "System, dnasr281@gmail.com"
    14 = "System, dnasr281@gmail.com"
    15 = "System, dnasr281@gmail.com"
    17 = "System, dnasr281@gmail.com"
    18 = "System, dnasr281@gmail.com"
    19 = "System, dnasr281@gmail.com"
    20 = "System, dnasr281@gmail.com"
    21 = "System, dnasr281@gmail.com"
    22 = "System, dnasr281@gmail.com"
    24 = "System, dnasr281@gmail.com"
    29 = "system, backup@backdoor.com, System"
    30 = "System, dnasr281@gmail.com"
    31 = "backup@backdoor.com, System"
    32 = "backup@backdoor.com, System"
    33 = "System, dnasr281@gmail.com"
    35 = "backup@backdoor.com, System"
    37 = "System, dnasr281@gmail.com"
    38 = "System, dnasr281@gmail.com"
    39 = "System, dnasr281@gmail.com"
    40 = "System, dnasr281@gmail.com"
    41 = "System, dnasr281@gmail.com"
    42 = "System, dnasr281@gmail.com"
    44 = "System, dnasr281@gmail.com"
    45 = "System, dnasr281@gmail.com"
    46 = "System, dnasr281@gmail.com"
    47 = "System, dnasr281@gmail.com"
    48 = "System, dnasr281@gmail.com"
    49 = "System, dnasr281@gmail.com"
    51 = "System, dnasr281@gmail.com"
    56 = "system, backup@backdoor.com, System"
    57 = "System, dnasr281@gmail.com"
    58 = "backup@backdoor.com, System"
    59 = "backup@backdoor.com, System"
    60 = "System, dnasr281@gmail.com"
    62 = "backup@backdoor.com, System"
    64 = "System, dnasr281@gmail.com"
    65 = "System, dnasr281@gmail.com"
    66 = "System, dnasr281@gmail.com"
    67 = "System, dnasr281@gmail.com"
    68 = "System, dnasr281@gmail.com"
    69 = "System, dnasr281@gmail.com"
    71 = "System, dnasr281@gmail.com"
    72 = "System, dnasr281@gmail.com"
    73 = "System, dnasr281@gmail.com"
    74 = "System, dnasr281@gmail.com"
    75 = "System, dnasr281@gmail.com"
    76 = "System, dnasr281@gmail.com"
    78 = "System, dnasr281@gmail.com"
    83 = "backup@backdoor.com, System"
    84 = "backup@backdoor.com, System"
    85 = "backup@backdoor.com, System"
    86 = "System, dnasr281@gmail.com"
    87 = "System, dnasr281@gmail.com"
    88 = "System, dnasr281@gmail.com"
    89 = "System, dnasr281@gmail.com"
    90 = "admin@admin.com, dnasr281@gmail.com"
    93 = "System, dnasr281@gmail.com"
    95 = "System, dnasr281@gmail.com"
    96 = "System, dnasr281@gmail.com"
    97 = "System, dnasr281@gmail.com"
    99 = "System, dnasr281@gmail.com"
    102 = "System, dnasr281@gmail.com"
    104 = "System, dnasr281@gmail.com"
    109 = "backup@backdoor.com, System"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    112 = "System, dnasr281@gmail.com"
    113 = "System, dnasr281@gmail.com"
    114 = "System, dnasr281@gmail.com"
    115 = "System, dnasr281@gmail.com"
    116 = "admin@admin.com, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    121 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    123 = "System, dnasr281@gmail.com"
    125 = "System, dnasr281@gmail.com"
    128 = "System, dnasr281@gmail.com"
    130 = "System, dnasr281@gmail.com"
    135 = "backup@backdoor.com, System"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    138 = "System, dnasr281@gmail.com"
    139 = "System, dnasr281@gmail.com"
    140 = "System, dnasr281@gmail.com"
    141 = "System, dnasr281@gmail.com"
    142 = "admin@admin.com, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
    147 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
    149 = "System, dnasr281@gmail.com"
    151 = "System, dnasr281@gmail.com"
    154 = "System, dnasr281@gmail.com"
    156 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
